# Weekly update: insert two new price records (Betarraga, Feria Lagunitas
# de Puerto Montt) ahead of the existing historical rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 293, pushing the existing rows 293:311 down
# to 295:313 (and extending the sheet dimension to A1:R313).
$ws.Rows.Item(293).Insert()
$ws.Rows.Item(293).Insert()

# New row 293: Primera grade, new 2022-07-04 (serial 44746) reading.
$ws.Range("A293").Value = 4
$ws.Range("B293").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C293").Value = "Los Lagos"
$ws.Range("D293").Value = 44746
$ws.Range("E293").Value = 10
$ws.Range("F293").Value = 100114014
$ws.Range("G293").Value = "Betarraga"
$ws.Range("H293").Value = "Sin especificar"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 250
$ws.Range("K293").Value = 1000
$ws.Range("L293").Value = 1000
$ws.Range("M293").Value = 1000
$ws.Range("N293").Value = "`$/paquete 5 unidades"
$ws.Range("O293").Value = "Región del Maule"
$ws.Range("P293").Value = 200
$ws.Range("Q293").Value = 5
$ws.Range("R293").Value = "Hortaliza"

# New row 294: Segunda grade, same date, same market.
$ws.Range("A294").Value = 4
$ws.Range("B294").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C294").Value = "Los Lagos"
$ws.Range("D294").Value = 44746
$ws.Range("E294").Value = 10
$ws.Range("F294").Value = 100114014
$ws.Range("G294").Value = "Betarraga"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Segunda"
$ws.Range("J294").Value = 250
$ws.Range("K294").Value = 1200
$ws.Range("L294").Value = 1200
$ws.Range("M294").Value = 1200
$ws.Range("N294").Value = "`$/paquete 5 unidades"
$ws.Range("O294").Value = "Región del Maule"
$ws.Range("P294").Value = 240
$ws.Range("Q294").Value = 5
$ws.Range("R294").Value = "Hortaliza"
